# Rename worksheet "Data_Final" to "Data-Final" to avoid loading errors
# when the workbook is read from Jupyter (e.g. pandas/openpyxl users
# treating underscores specially).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data_Final")
$ws.Name = "Data-Final"
